$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) "step 3 character dutch;" bullet gets the green author color
#    (00A933) that's used throughout the rest of the list, on both the
#    paragraph mark and the run.
# ---------------------------------------------------------------------
$rng1 = $d.Content
$rng1.Find.ClearFormatting()
$found1 = $rng1.Find.Execute("step 3 character dutch;")
if (-not $found1) { throw "could not find 'step 3 character dutch;'" }
$para1 = $rng1.Paragraphs(1).Range
$para1.Font.Color = 3385600   # wdColor for 00A933 (R=0,G=169,B=51)

# ---------------------------------------------------------------------
# 2) Merge the three runs of the CSV-delimiter bullet into a single run
#    with the combined text (no more mid-sentence run splits).
# ---------------------------------------------------------------------
$rngCsv = $d.Content
$rngCsv.Find.ClearFormatting()
$rngCsv.Find.Replacement.ClearFormatting()
$rngCsv.Find.Execute(
    "add CSV delimiter and quote type in first step or auto-detect;",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "add CSV delimiter and quote type in first step or auto-detect;", 2)

# ---------------------------------------------------------------------
# 3) Merge the "s" + "elect localization just for " runs (which precede
#    the differently-fonted "CSV" run) into a single run reading
#    "select localization just for ".
# ---------------------------------------------------------------------
$rngSel = $d.Content
$rngSel.Find.ClearFormatting()
$rngSel.Find.Replacement.ClearFormatting()
$rngSel.Find.Execute(
    "select localization just for ",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "select localization just for ", 2)

# ---------------------------------------------------------------------
# 4) Normal style: turn overflow punctuation on
#    (w:overflowPunct false -> true), i.e.
#    ParagraphFormat.HangingPunctuation = True.
# ---------------------------------------------------------------------
$normal = $d.Styles.Item("Normal")
$normal.ParagraphFormat.HangingPunctuation = $true

Write-Output "done"
